$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D and E (batsman column shifts from D to F)
$ws.Range("D1:E1").EntireColumn.Insert()

# Header row
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"

# Row 2 data
$ws.Range("D2").Value = "Delhi Capitals"
$ws.Range("E2").Value = "Mumbai Indians"

# Row 3 data
$ws.Range("D3").Value = "Delhi Capitals"
$ws.Range("E3").Value = "Rajasthan Royals"
